$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared string "sCs" is introduced alongside ECs/FAPs as a third
# sending/target cluster, and rows 4-10 are added to give the full 3x3
# sending-cluster x target-cluster matrix for the Ntf3->Ntrk1 pair.

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ntf3"
$ws.Cells.Item(2, 3).Value = "Ntrk1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 8.859944
$ws.Cells.Item(2, 8).Value = 26.579832
$ws.Cells.Item(2, 9).Value = 0.4921128329655918
$ws.Cells.Item(2, 10).Value = 0.4921128329655918
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.0002123333333333333
$ws.Cells.Item(2, 14).Value = 0.000637
$ws.Cells.Item(2, 15).Value = 0.0008844802186068535
$ws.Cells.Item(2, 16).Value = 0.0008844802186068534
$ws.Cells.Item(2, 17).Value = 0.001881261442666667
$ws.Cells.Item(2, 18).Value = 0.016931352984
$ws.Cells.Item(2, 19).Value = 0.0004352640660806446
$ws.Cells.Item(2, 20).Value = 0.0004352640660806446

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ntf3"
$ws.Cells.Item(3, 3).Value = "Ntrk1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 8.859944
$ws.Cells.Item(3, 8).Value = 26.579832
$ws.Cells.Item(3, 9).Value = 0.4921128329655918
$ws.Cells.Item(3, 10).Value = 0.4921128329655918
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.220618
$ws.Cells.Item(3, 14).Value = 0.6618539999999999
$ws.Cells.Item(3, 15).Value = 0.9189902207312721
$ws.Cells.Item(3, 16).Value = 0.9189902207312721
$ws.Cells.Item(3, 17).Value = 1.954663125392
$ws.Cells.Item(3, 18).Value = 17.591968128528
$ws.Cells.Item(3, 19).Value = 0.4522468809917408
$ws.Cells.Item(3, 20).Value = 0.4522468809917408

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Ntf3"
$ws.Cells.Item(4, 3).Value = "Ntrk1"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 8.859944
$ws.Cells.Item(4, 8).Value = 26.579832
$ws.Cells.Item(4, 9).Value = 0.4921128329655918
$ws.Cells.Item(4, 10).Value = 0.4921128329655918
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.01923533333333333
$ws.Cells.Item(4, 14).Value = 0.057706
$ws.Cells.Item(4, 15).Value = 0.08012529905012102
$ws.Cells.Item(4, 16).Value = 0.08012529905012102
$ws.Cells.Item(4, 17).Value = 0.1704239761546667
$ws.Cells.Item(4, 18).Value = 1.533815785392
$ws.Cells.Item(4, 19).Value = 0.0394306879077703
$ws.Cells.Item(4, 20).Value = 0.0394306879077703

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Ntf3"
$ws.Cells.Item(5, 3).Value = "Ntrk1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 5.867977666666667
$ws.Cells.Item(5, 8).Value = 17.603933
$ws.Cells.Item(5, 9).Value = 0.3259283708025871
$ws.Cells.Item(5, 10).Value = 0.3259283708025871
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.0002123333333333333
$ws.Cells.Item(5, 14).Value = 0.000637
$ws.Cells.Item(5, 15).Value = 0.0008844802186068535
$ws.Cells.Item(5, 16).Value = 0.0008844802186068534
$ws.Cells.Item(5, 17).Value = 0.001245967257888889
$ws.Cells.Item(5, 18).Value = 0.011213705321
$ws.Cells.Item(5, 19).Value = 0.0002882771966576478
$ws.Cells.Item(5, 20).Value = 0.0002882771966576478

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Ntf3"
$ws.Cells.Item(6, 3).Value = "Ntrk1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 5.867977666666667
$ws.Cells.Item(6, 8).Value = 17.603933
$ws.Cells.Item(6, 9).Value = 0.3259283708025871
$ws.Cells.Item(6, 10).Value = 0.3259283708025871
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.220618
$ws.Cells.Item(6, 14).Value = 0.6618539999999999
$ws.Cells.Item(6, 15).Value = 0.9189902207312721
$ws.Cells.Item(6, 16).Value = 0.9189902207312721
$ws.Cells.Item(6, 17).Value = 1.294581496864667
$ws.Cells.Item(6, 18).Value = 11.651233471782
$ws.Cells.Item(6, 19).Value = 0.2995249854264534
$ws.Cells.Item(6, 20).Value = 0.2995249854264534

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Ntf3"
$ws.Cells.Item(7, 3).Value = "Ntrk1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 5.867977666666667
$ws.Cells.Item(7, 8).Value = 17.603933
$ws.Cells.Item(7, 9).Value = 0.3259283708025871
$ws.Cells.Item(7, 10).Value = 0.3259283708025871
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.01923533333333333
$ws.Cells.Item(7, 14).Value = 0.057706
$ws.Cells.Item(7, 15).Value = 0.08012529905012102
$ws.Cells.Item(7, 16).Value = 0.08012529905012102
$ws.Cells.Item(7, 17).Value = 0.1128725064108889
$ws.Cells.Item(7, 18).Value = 1.015852557698
$ws.Cells.Item(7, 19).Value = 0.02611510817947602
$ws.Cells.Item(7, 20).Value = 0.02611510817947602

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Ntf3"
$ws.Cells.Item(8, 3).Value = "Ntrk1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 3.275965666666667
$ws.Cells.Item(8, 8).Value = 9.827897
$ws.Cells.Item(8, 9).Value = 0.1819587962318212
$ws.Cells.Item(8, 10).Value = 0.1819587962318212
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.0002123333333333333
$ws.Cells.Item(8, 14).Value = 0.000637
$ws.Cells.Item(8, 15).Value = 0.0008844802186068535
$ws.Cells.Item(8, 16).Value = 0.0008844802186068534
$ws.Cells.Item(8, 17).Value = 0.0006955967098888889
$ws.Cells.Item(8, 18).Value = 0.006260370389
$ws.Cells.Item(8, 19).Value = 0.0001609389558685611
$ws.Cells.Item(8, 20).Value = 0.0001609389558685611

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Ntf3"
$ws.Cells.Item(9, 3).Value = "Ntrk1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 3.275965666666667
$ws.Cells.Item(9, 8).Value = 9.827897
$ws.Cells.Item(9, 9).Value = 0.1819587962318212
$ws.Cells.Item(9, 10).Value = 0.1819587962318212
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.220618
$ws.Cells.Item(9, 14).Value = 0.6618539999999999
$ws.Cells.Item(9, 15).Value = 0.9189902207312721
$ws.Cells.Item(9, 16).Value = 0.9189902207312721
$ws.Cells.Item(9, 17).Value = 0.7227369934486666
$ws.Cells.Item(9, 18).Value = 6.504632941037999
$ws.Cells.Item(9, 19).Value = 0.1672183543130779
$ws.Cells.Item(9, 20).Value = 0.1672183543130779

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Ntf3"
$ws.Cells.Item(10, 3).Value = "Ntrk1"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 3.275965666666667
$ws.Cells.Item(10, 8).Value = 9.827897
$ws.Cells.Item(10, 9).Value = 0.1819587962318212
$ws.Cells.Item(10, 10).Value = 0.1819587962318212
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.01923533333333333
$ws.Cells.Item(10, 14).Value = 0.057706
$ws.Cells.Item(10, 15).Value = 0.08012529905012102
$ws.Cells.Item(10, 16).Value = 0.08012529905012102
$ws.Cells.Item(10, 17).Value = 0.06301429158688888
$ws.Cells.Item(10, 18).Value = 0.5671286242820001
$ws.Cells.Item(10, 19).Value = 0.0145795029628747
$ws.Cells.Item(10, 20).Value = 0.0145795029628747

